$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "25.888.16"
$ws.Range("E2").Value = "  -0.91%  "

# Row 3
$ws.Range("D3").Value = "1.738.09"
$ws.Range("E3").Value = "  -1.81%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.04"
$ws.Range("E5").Value = "  -3.51%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9997"
$ws.Range("E6").Value = "  -0.10%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5252"
$ws.Range("E7").Value = "  +0.16%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2759"
$ws.Range("E8").Value = "  -0.26%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "39.39"
$ws.Range("E9").Value = "  -2.84%  "

# Row 10
$ws.Range("E10").Value = "  -1.05%  "

# Row 11
$ws.Range("D11").Value = "1.738.34"
$ws.Range("E11").Value = "  -2.12%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07099"
$ws.Range("E12").Value = "  +1.11%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "15.13"
$ws.Range("E13").Value = "  -4.49%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6424"
$ws.Range("E14").Value = "  -0.80%  "

# Row 15
$ws.Range("E15").Value = "  -0.35%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "76.96"
$ws.Range("E16").Value = "  -1.79%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9998"
$ws.Range("E17").Value = "  -0.14%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.9993"
$ws.Range("E18").Value = "  -0.11%  "

# Row 19
$ws.Range("D19").Value = "25.864.37"
$ws.Range("E19").Value = "  -1.07%  "

# Row 20
$ws.Range("E20").Value = "  -1.51%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.000006667"
$ws.Range("E21").Value = "  -1.80%  "

# Row 22
$ws.Range("D22").Value = "1.960.28"
$ws.Range("E22").Value = "  -1.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.275"
$ws.Range("E23").Value = "  +4.62%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "8.811"
$ws.Range("E24").Value = "  +3.61%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.167"
$ws.Range("E25").Value = "  -0.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "140.30"
$ws.Range("E26").Value = "  +0.87%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "1.513"
$ws.Range("E27").Value = "  +1.24%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.17"
$ws.Range("E28").Value = "  -0.33%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.789"
$ws.Range("E29").Value = "  -3.98%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "102.37"
$ws.Range("E30").Value = "  -0.69%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08331"
$ws.Range("E31").Value = "  -1.14%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.726"
$ws.Range("E32").Value = "  +0.15%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.596"
$ws.Range("E33").Value = "  +3.47%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.04511"
$ws.Range("E34").Value = "  +0.54%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.609"
$ws.Range("E35").Value = "  -1.74%  "

# Row 36
$ws.Range("B36").Value = "ARBITRUM"
$ws.Range("C36").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.9726"
$ws.Range("E36").Value = "  -3.92%  "

# Row 37
$ws.Range("B37").Value = "ImmutableX"
$ws.Range("C37").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.6196"
$ws.Range("E37").Value = "  +1.98%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.684"
$ws.Range("E38").Value = "  -2.58%  "

# Row 39
$ws.Range("E39").Value = "  -1.10%  "

# Row 40
$ws.Range("B40").Value = "PaxDollar"
$ws.Range("C40").Value = "https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.9992"
$ws.Range("E40").Value = "  -0.28%  "

# Row 41
$ws.Range("B41").Value = "RenderToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.902"
$ws.Range("E41").Value = "  -4.81%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "99.88"
$ws.Range("E42").Value = "  -3.22%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.3852"
$ws.Range("E43").Value = "  -1.13%  "

# Row 44
$ws.Range("B44").Value = "TrustWalletToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.7272"
$ws.Range("E44").Value = "  -3.46%  "

# Row 45
$ws.Range("B45").Value = "FraxShare"
$ws.Range("C45").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.014"
$ws.Range("E45").Value = "  +1.39%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.05324"
$ws.Range("E46").Value = "  -3.59%  "

# Row 47
$ws.Range("E47").Value = "  -0.21%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.213"
$ws.Range("E48").Value = "  -3.01%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "53.62"
$ws.Range("E49").Value = "  +1.60%  "

# Row 50
$ws.Range("E50").Value = "  -1.02%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.638"
$ws.Range("E51").Value = "  +1.49%  "
